$d = $word.ActiveDocument

# 1) Reword the "Fue desarrollada..." sentence start.
$d.Content.Find.Execute(
    "Fue desarrollada con un enfoque en permitir",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "Esta biblioteca fue desarrollada con el enfoque de permitir",
    2)

# 2) Append new paragraphs at the very end of the document (after the last
#    paragraph, before the final section break).
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()

$newParas = @(
    "En este caso, se llevó a cabo la construcción de un modelo de Red Neuronal Convolucional (CNN, por sus siglas en inglés) para la detección de peatones utilizando imágenes. El uso de CNNs se ha vuelto muy popular en la detección de objetos en imágenes, gracias a su capacidad para extraer características importantes a través de diferentes niveles de abstracción.",
    "El proceso inició con el preprocesamiento de las imágenes, que implicó la lectura, redimensionamiento y normalización de las imágenes. Esta etapa es crucial ya que prepara los datos para la alimentación de la red neuronal. Posteriormente, las imágenes preprocesadas fueron divididas en dos categorías, una etiquetada como 'peatón' y la otra como 'no peatón', creando así un conjunto de datos etiquetados que facilita el entrenamiento supervisado.",
    'El diseño de la CNN incluyó varias capas convolucionales intercaladas con capas de agrupamiento o "pooling". Las capas convolucionales son responsables de la detección de características visuales en las imágenes, mientras que las capas de pooling reducen la dimensionalidad de los datos, lo que ahorra tiempo de cálculo y ayuda a evitar el sobreajuste. ',
    'Después de las capas de convolución y pooling, la entrada se aplanó para poder ser alimentada a las capas densas, que son responsables de la clasificación de las características aprendidas en las capas previas. Se incorporó una capa de abandono o "dropout" después de las capas densas para mejorar la generalización del modelo y minimizar el sobreajuste. Finalmente, la capa final, una capa densa con una sola neurona y una función de activación sigmoide, generó la salida del modelo, proporcionando la probabilidad de que la imagen contenga un peatón.',
    "El modelo se compiló con una función de pérdida de entropía cruzada binaria, dado que el problema de detección de peatones es un problema de clasificación binaria, y se utilizó el optimizador 'rmsprop'. ",
    "Se realizó una validación cruzada estratificada de 5 divisiones para evaluar la eficacia del modelo, proporcionando una visión más precisa de cómo se desempeñará el modelo en datos no vistos. ",
    "Los resultados de la validación cruzada, representados en forma de gráficos de precisión y pérdida para cada pliegue, ofrecieron un medio visual para evaluar el rendimiento del modelo a lo largo de las épocas, tanto en el conjunto de entrenamiento como en el de validación. Al final, se calculó la precisión promedio y la desviación estándar de los 5 pliegues, proporcionando una medida general del rendimiento del modelo.",
    "En resumen, la CNN demostró ser una herramienta eficaz para la detección de peatones. El uso de la validación cruzada proporcionó una evaluación robusta del modelo, y la inclusión de capas de dropout y pooling ayudó a minimizar el sobreajuste."
)

$last = $d.Paragraphs.Last
for ($i = 0; $i -lt $newParas.Length; $i++) {
    $last.Range.Text = $newParas[$i]
    $last.Range.Font.Reset()
    if ($i -lt ($newParas.Length - 1)) {
        $last.Range.InsertParagraphAfter()
        $last = $d.Paragraphs.Last
    }
}
